$wb = $excel.ActiveWorkbook

# Sheet1: Outbreak_Locations - add row 38
$ws1 = $wb.Worksheets.Item("Outbreak_Locations")
$ws1.Range("A38").Value = 16
$ws1.Range("B38").Value = 650
$ws1.Range("C38").Value = 750

# Sheet2: Store_Locations - add row 39
$ws2 = $wb.Worksheets.Item("Store_Locations")
$ws2.Range("A39").Value = 16
$ws2.Range("B39").Value = 550
$ws2.Range("C39").Value = 550
$ws2.Range("D39").Value = "Chain 1"

# Sheet3: Population - add row 17 (lower bound row)
$ws3 = $wb.Worksheets.Item("Population")

# Copy formatting: B/C from row 16 (style s=1, has value), D/E from row 2 (style
# s=2, empty cell formatting) to match the "empty" cells in the new row.
$ws3.Range("B16:C16").Copy()
$ws3.Range("B17:C17").PasteSpecial(-4122)
$ws3.Range("D2:E2").Copy()
$ws3.Range("D17:E17").PasteSpecial(-4122)

# For A17: base on the D2 (s=2) formatting (applyFill + applyBorder already set),
# then strip fill and change border shape to left/right-only thin border.
$ws3.Range("D2").Copy()
$ws3.Range("A17").PasteSpecial(-4122)
$ws3.Range("A17").Interior.Pattern = -4142
$ws3.Range("A17").Borders.Item(8).LineStyle = -4142
$ws3.Range("A17").Borders.Item(9).LineStyle = -4142
$ws3.Range("A17").Borders.Item(7).LineStyle = 1
$ws3.Range("A17").Borders.Item(7).Weight = 2
$ws3.Range("A17").Borders.Item(10).LineStyle = 1
$ws3.Range("A17").Borders.Item(10).Weight = 2

# Now set the values
$ws3.Range("A17").Value = 16
$ws3.Range("B17").Value = "uniform"
$ws3.Range("C17").Value = 500
$ws3.Range("D17").ClearContents()
$ws3.Range("E17").ClearContents()
